# Update cryptocurrency price/volume data (and fix a NEARProtocol/Maker
# row ordering swap) per the Fri Feb 16 11:52:39 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''52.247.15'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '''2.825.05'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''356.84'
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").Value = '''112.24'
$ws.Range("E6").Value = '  -3.18%  '
$ws.Range("D7").Value = '''0.572'
$ws.Range("E7").Value = '  +4.07%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '''0.600'
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").Value = '''41.05'
$ws.Range("E10").Value = '  -4.10%  '
$ws.Range("D11").Value = '''0.0863'
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = '''19.92'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = '''7.78'
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").Value = '''3.274.75'
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("D16").Value = '''2.838.56'
$ws.Range("E16").Value = '  +2.43%  '
$ws.Range("E17").Value = '  +4.47%  '
$ws.Range("D18").Value = '''52.218.35'
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = '''7.50'
$ws.Range("E19").Value = '  +4.03%  '
$ws.Range("D20").Value = '''3.19'
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("D22").Value = '''0.0₃0996'
$ws.Range("E22").Value = '  +1.78%  '
$ws.Range("D23").Value = '''70.54'
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").Value = '''271.03'
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("E25").Value = '  +1.76%  '
$ws.Range("D26").Value = '''26.98'
$ws.Range("E26").Value = '  +1.36%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = '''10.34'
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("D30").Value = '''0.0487'
$ws.Range("E30").Value = '  +8.88%  '
$ws.Range("D31").Value = '''0.143'
$ws.Range("E31").Value = '  +2.62%  '
$ws.Range("D32").Value = '''52.52'
$ws.Range("E32").Value = '  +4.66%  '
$ws.Range("D33").Value = '''35.23'
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("D34").Value = '''5.94'
$ws.Range("E34").Value = '  +4.20%  '
$ws.Range("D35").Value = '''5.59'
$ws.Range("E35").Value = '  +12.64%  '
$ws.Range("D36").Value = '''0.0856'
$ws.Range("E36").Value = '  +3.78%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("D39").Value = '''2.04'
$ws.Range("E39").Value = '  -3.14%  '
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("D41").Value = '''0.117'
$ws.Range("E41").Value = '  +1.90%  '
$ws.Range("D42").Value = '''127.53'
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("E43").Value = '  -3.98%  '
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '''2.092.46'
$ws.Range("E46").Value = '  +1.41%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''3.37'
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("D49").Value = '''5.93'
$ws.Range("E49").Value = '  +7.32%  '
$ws.Range("D50").Value = '''0.965'
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("D51").Value = '''9.17'
$ws.Range("E51").Value = '  +3.10%  '
